# Include the image title (e.g. from markdown `![alt text](link "title")`)
# in PowerPoint's description of the image, alongside the link, for images
# that are "implicit figures" (i.e. a standalone image with a caption
# underneath it). Previously the title was ignored when writing the pptx;
# now it is prepended to the existing link-only description as
# "<title>  <link>" (here the sentinel title used for implicit figures is
# "fig:").

$p = $ppt.ActivePresentation

# Slide 2: picture ("lalune.jpg") followed by a "The Moon" caption textbox
# -> this is an implicit figure, so its description gains the title prefix.
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).AlternativeText = "fig:  lalune.jpg"

# Slide 3: title, then picture ("lalune.jpg"), then a "The Moon" caption
# textbox -> also an implicit figure.
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).AlternativeText = "fig:  lalune.jpg"
